$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose match data (everything except the running index in column A,
# and the constant Div/Div Original Name/Date columns C:E) got swapped between
# the two rows during this update.
$pairs = @(
    @(7, 8),
    @(32, 33),
    @(67, 68),
    @(77, 78),
    @(179, 180)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Column B (id)
    $b1 = $ws.Cells.Item($r1, 2).Value2
    $b2 = $ws.Cells.Item($r2, 2).Value2
    $ws.Cells.Item($r1, 2).Value2 = $b2
    $ws.Cells.Item($r2, 2).Value2 = $b1

    # Columns F through AC (6 through 29): HomeTeam, AwayTeam, FTHG, FTAG,
    # FTR, and all odds/PL columns.
    for ($col = 6; $col -le 29; $col++) {
        $v1 = $ws.Cells.Item($r1, $col).Value2
        $v2 = $ws.Cells.Item($r2, $col).Value2
        $ws.Cells.Item($r1, $col).Value2 = $v2
        $ws.Cells.Item($r2, $col).Value2 = $v1
    }
}
